# Update Leve profit-tracking values (currentAveragePrice* / LevePrice* / LeveProfit*)
# across the Belias_Profits leve sheets, per the scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused
$ws.Range("H33").Value = 34483188
$ws.Range("I33").Value = 132.7619
$ws.Range("J33").Value = 125001210
$ws.Range("K33").Value = 132.7619
$ws.Range("L33").Value = 125001210
$ws.Range("M33").Value = 96.2381
$ws.Range("N33").Value = -125001668
# Row 49: Going Nowhere Fast
$ws.Range("H49").Value = 1433.2858
$ws.Range("I49").Value = 886.8
$ws.Range("J49").Value = 2799.5
$ws.Range("K49").Value = 2660.4
$ws.Range("L49").Value = 8398.5
$ws.Range("M49").Value = -2524.4
$ws.Range("N49").Value = -8670.5
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 3083.8413
$ws.Range("I132").Value = 1505.6938
$ws.Range("K132").Value = 4517.0814
$ws.Range("M132").Value = -1987.0814
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1356492.9
$ws.Range("I137").Value = 1463.0278
$ws.Range("J137").Value = 11112708
$ws.Range("K137").Value = 4389.0834
$ws.Range("L137").Value = 33338124
$ws.Range("M137").Value = -1839.0834
$ws.Range("N137").Value = -33343224
# Row 138: All-night Crafting
$ws.Range("H138").Value = 3468408.8
$ws.Range("J138").Value = 7735766
$ws.Range("L138").Value = 23207298
$ws.Range("N138").Value = -23217578

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 3660.59
$ws.Range("I32").Value = 2788.5327
$ws.Range("J32").Value = 13689.25
$ws.Range("K32").Value = 2788.5327
$ws.Range("L32").Value = 13689.25
$ws.Range("M32").Value = -2501.5327
$ws.Range("N32").Value = -14263.25
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 1472.5883
$ws.Range("I61").Value = 1302.0571
$ws.Range("J61").Value = 1845.625
$ws.Range("K61").Value = 1302.0571
$ws.Range("L61").Value = 1845.625
$ws.Range("M61").Value = -1090.0571
$ws.Range("N61").Value = -2269.625
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 38347.234
$ws.Range("I74").Value = 42933.293
$ws.Range("J74").Value = 20003
$ws.Range("K74").Value = 42933.293
$ws.Range("L74").Value = 20003
$ws.Range("M74").Value = -42059.293
$ws.Range("N74").Value = -21751
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 38347.234
$ws.Range("I77").Value = 42933.293
$ws.Range("J77").Value = 20003
$ws.Range("K77").Value = 214666.465
$ws.Range("L77").Value = 100015
$ws.Range("M77").Value = -210298.465
$ws.Range("N77").Value = -108751
# Row 101: Art Imitates Life
$ws.Range("H101").Value = 29750
$ws.Range("J101").Value = 29750
$ws.Range("L101").Value = 29750
$ws.Range("N101").Value = -36240
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 1472.5883
$ws.Range("I136").Value = 1302.0571
$ws.Range("J136").Value = 1845.625
$ws.Range("K136").Value = 3906.1713
$ws.Range("L136").Value = 5536.875
$ws.Range("M136").Value = -1356.1713
$ws.Range("N136").Value = -10636.875

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 731598.8
$ws.Range("I134").Value = 1254203.5
$ws.Range("J134").Value = 4496.609
$ws.Range("K134").Value = 3762610.5
$ws.Range("L134").Value = 13489.827
$ws.Range("M134").Value = -3760075.5
$ws.Range("N134").Value = -18559.827

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 611.4583
$ws.Range("I22").Value = 185
$ws.Range("J22").Value = 916.0714
$ws.Range("K22").Value = 185
$ws.Range("L22").Value = 916.0714
$ws.Range("M22").Value = 165
$ws.Range("N22").Value = -1616.0714
# Row 31: Wall Not Found
$ws.Range("H31").Value = 13890176
$ws.Range("I31").Value = 1191.4375
$ws.Range("J31").Value = 125002056
$ws.Range("K31").Value = 1191.4375
$ws.Range("L31").Value = 125002056
$ws.Range("M31").Value = -896.4375
$ws.Range("N31").Value = -125002646
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 13890176
$ws.Range("I34").Value = 1191.4375
$ws.Range("J34").Value = 125002056
$ws.Range("K34").Value = 1191.4375
$ws.Range("L34").Value = 125002056
$ws.Range("M34").Value = -989.4375
$ws.Range("N34").Value = -125002460
# Row 102: The Ear Is the Way to the Heart
$ws.Range("H102").Value = 21432.5
$ws.Range("J102").Value = 21432.5
$ws.Range("L102").Value = 21432.5
$ws.Range("N102").Value = -26300.5
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 4404.972
$ws.Range("I134").Value = 4053.9092
$ws.Range("J134").Value = 8266.666999999999
$ws.Range("K134").Value = 12161.7276
$ws.Range("L134").Value = 24800.001
$ws.Range("M134").Value = -9626.7276
$ws.Range("N134").Value = -29870.001

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water
$ws.Range("H4").Value = 899.8946999999999
$ws.Range("I4").Value = 293.5
$ws.Range("J4").Value = 1340.909
$ws.Range("K4").Value = 880.5
$ws.Range("L4").Value = 4022.727
$ws.Range("M4").Value = -768.5
$ws.Range("N4").Value = -4246.727000000001
# Row 5: What a Sap
$ws.Range("H5").Value = 809.16
$ws.Range("I5").Value = 525.6429000000001
$ws.Range("J5").Value = 1170
$ws.Range("K5").Value = 1576.9287
$ws.Range("L5").Value = 3510
$ws.Range("M5").Value = -1464.9287
$ws.Range("N5").Value = -3734
# Row 17: Chew the Fat
$ws.Range("H17").Value = 300.25
$ws.Range("I17").Value = 199.5
$ws.Range("J17").Value = 401
$ws.Range("K17").Value = 598.5
$ws.Range("L17").Value = 1203
$ws.Range("M17").Value = -429.5
$ws.Range("N17").Value = -1541
# Row 33: Cooking with Gas
$ws.Range("H33").Value = 68.333336
$ws.Range("I33").Value = 61.666668
$ws.Range("J33").Value = 75
$ws.Range("K33").Value = 370.000008
$ws.Range("L33").Value = 450
$ws.Range("M33").Value = -87.00000799999998
$ws.Range("N33").Value = -1016
# Row 58: Bread in the Clouds
$ws.Range("H58").Value = 3233.0952
$ws.Range("I58").Value = 2002.5
$ws.Range("J58").Value = 3362.6316
$ws.Range("K58").Value = 6007.5
$ws.Range("L58").Value = 10087.8948
$ws.Range("M58").Value = -5879.5
$ws.Range("N58").Value = -10343.8948
# Row 61: Red Letter Day
$ws.Range("H61").Value = 456.33334
$ws.Range("I61").Value = 257.7143
$ws.Range("J61").Value = 582.7273
$ws.Range("K61").Value = 773.1428999999999
$ws.Range("L61").Value = 1748.1819
$ws.Range("M61").Value = -558.1428999999999
$ws.Range("N61").Value = -2178.1819
# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 645179.8
$ws.Range("I113").Value = 1122766.1
$ws.Range("J113").Value = 438.4
$ws.Range("K113").Value = 3368298.3
$ws.Range("L113").Value = 1315.2
$ws.Range("M113").Value = -3366128.3
$ws.Range("N113").Value = -5655.2
# Row 122: Salt of the North
$ws.Range("H122").Value = 43000.92
$ws.Range("I122").Value = 55511.25
$ws.Range("J122").Value = 1299.8334
$ws.Range("K122").Value = 499601.25
$ws.Range("L122").Value = 11698.5006
$ws.Range("M122").Value = -497151.25
$ws.Range("N122").Value = -16598.5006
# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 930.75
$ws.Range("I131").Value = 840
$ws.Range("J131").Value = 934.53125
$ws.Range("K131").Value = 2520
$ws.Range("L131").Value = 2803.59375
$ws.Range("M131").Value = 2520
$ws.Range("N131").Value = -12883.59375
# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 809.16
$ws.Range("I135").Value = 525.6429000000001
$ws.Range("J135").Value = 1170
$ws.Range("K135").Value = 4730.7861
$ws.Range("L135").Value = 10530
$ws.Range("M135").Value = -2195.7861
$ws.Range("N135").Value = -15600

$ws = $wb.Worksheets.Item("GSM")
# Row 5: Hora at Me
$ws.Range("H5").Value = 10680.833
$ws.Range("I5").Value = 5800
$ws.Range("J5").Value = 13121.25
$ws.Range("K5").Value = 5800
$ws.Range("L5").Value = 13121.25
$ws.Range("M5").Value = -5688
$ws.Range("N5").Value = -13345.25
# Row 64: Halonic Hermeneutics
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
# Row 67: Transposing Theology (L)
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
# Row 132: On Board for Lar
$ws.Range("H132").Value = 1251821.8
$ws.Range("I132").Value = 1555.3137
$ws.Range("J132").Value = 3450566.2
$ws.Range("K132").Value = 4665.9411
$ws.Range("L132").Value = 10351698.6
$ws.Range("M132").Value = -2135.9411
$ws.Range("N132").Value = -10356758.6

$ws = $wb.Worksheets.Item("LTW")
# Row 60: Tenderfoot Moments
$ws.Range("H60").Value = 29000
$ws.Range("J60").Value = 29000
$ws.Range("L60").Value = 29000
$ws.Range("N60").Value = -30018
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 2850.4265
$ws.Range("I132").Value = 2469.1372
$ws.Range("K132").Value = 7407.4116
$ws.Range("M132").Value = -4877.4116

$ws = $wb.Worksheets.Item("WVR")
# Row 103: To the Tops
$ws.Range("H103").Value = 30200.334
$ws.Range("J103").Value = 30200.334
$ws.Range("L103").Value = 30200.334
$ws.Range("N103").Value = -32544.334
# Row 107: Flax Wax
$ws.Range("H107").Value = 9394.083000000001
$ws.Range("I107").Value = 15461.571
$ws.Range("J107").Value = 899.6
$ws.Range("K107").Value = 46384.713
$ws.Range("L107").Value = 2698.8
$ws.Range("M107").Value = -44464.713
$ws.Range("N107").Value = -6538.8
